# Auto-generated script applying cell updates from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "19.942.64"
$ws.Range("E2").Value = "  -7.77%  "
Set-TextValue "D3" "1.406.67"
$ws.Range("E3").Value = "  -8.30%  "
Set-TextValue "D4" "0.9976"
$ws.Range("E4").Value = "  -0.28%  "
Set-TextValue "D5" "0.9943"
$ws.Range("E5").Value = "  -0.55%  "
Set-TextValue "D6" "271.42"
$ws.Range("E6").Value = "  -6.25%  "
Set-TextValue "D7" "0.3671"
$ws.Range("E7").Value = "  -6.84%  "
Set-TextValue "D8" "0.3117"
$ws.Range("E8").Value = "  -1.64%  "
Set-TextValue "D9" "39.51"
$ws.Range("E9").Value = "  -6.67%  "
$ws.Range("E10").Value = "  -5.48%  "
Set-TextValue "D11" "0.06493"
$ws.Range("E11").Value = "  -9.43%  "
Set-TextValue "D12" "0.9988"
$ws.Range("E12").Value = "  -0.21%  "
Set-TextValue "D13" "5.407"
$ws.Range("E13").Value = "  -5.21%  "
Set-TextValue "D14" "17.36"
$ws.Range("E14").Value = "  -5.19%  "
Set-TextValue "D15" "6.144"
$ws.Range("E15").Value = "  -7.18%  "
Set-TextValue "D16" "1.402.85"
$ws.Range("E16").Value = "  -8.80%  "
Set-TextValue "D17" "0.00001017"
$ws.Range("E17").Value = "  -6.93%  "
Set-TextValue "D18" "0.05701"
$ws.Range("E18").Value = "  -13.64%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D19" "71.03"
$ws.Range("E19").Value = "  -15.37%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "0.9961"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "5.573"
$ws.Range("E21").Value = "  -9.35%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D22" "14.82"
$ws.Range("E22").Value = "  -4.55%  "
Set-TextValue "D23" "10.97"
$ws.Range("E23").Value = "  +2.29%  "
Set-TextValue "D24" "2.261"
$ws.Range("E24").Value = "  -4.61%  "
Set-TextValue "D25" "19.959.76"
$ws.Range("E25").Value = "  -7.68%  "
Set-TextValue "D26" "2.237"
$ws.Range("E26").Value = "  -4.91%  "
Set-TextValue "D27" "135.79"
$ws.Range("E27").Value = "  -9.70%  "
Set-TextValue "D28" "16.93"
$ws.Range("E28").Value = "  -7.65%  "
Set-TextValue "D29" "1.558.33"
$ws.Range("E29").Value = "  -8.75%  "
Set-TextValue "D30" "109.92"
$ws.Range("E30").Value = "  -6.37%  "
Set-TextValue "D31" "4.100"
$ws.Range("E31").Value = "  -15.46%  "
Set-TextValue "D32" "5.302"
$ws.Range("E32").Value = "  -12.61%  "
Set-TextValue "D33" "0.8287"
$ws.Range("E33").Value = "  -12.10%  "
Set-TextValue "D34" "0.07658"
$ws.Range("E34").Value = "  -5.96%  "
Set-TextValue "D35" "8.301"
$ws.Range("E35").Value = "  -2.08%  "
Set-TextValue "D36" "1.455"
$ws.Range("E36").Value = "  -0.37%  "
Set-TextValue "D37" "0.05785"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("E38").Value = "  -7.32%  "
Set-TextValue "D39" "0.9953"
$ws.Range("E39").Value = "  -0.37%  "
Set-TextValue "D40" "0.02071"
$ws.Range("E40").Value = "  -6.54%  "
Set-TextValue "D41" "0.1904"
$ws.Range("E41").Value = "  -6.41%  "
Set-TextValue "D42" "10.35"
$ws.Range("E42").Value = "  -6.60%  "
Set-TextValue "D43" "1.085"
$ws.Range("E43").Value = "  -8.44%  "
Set-TextValue "D44" "0.5283"
$ws.Range("E44").Value = "  -9.07%  "
Set-TextValue "D45" "12.25"
$ws.Range("E45").Value = "  -5.70%  "
Set-TextValue "D46" "3.491"
$ws.Range("E46").Value = "  -6.10%  "
Set-TextValue "D47" "0.5117"
$ws.Range("E47").Value = "  -7.90%  "
Set-TextValue "D48" "111.62"
$ws.Range("E48").Value = "  -3.98%  "
Set-TextValue "D49" "1.772"
$ws.Range("E49").Value = "  -5.81%  "
Set-TextValue "D50" "1.040"
$ws.Range("E50").Value = "  -10.89%  "
Set-TextValue "D51" "0.9977"
$ws.Range("E51").Value = "  -0.26%  "
